$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right=5, Wrong=-1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 (Total): Right=85, Wrong=-2.4, Max label updated accordingly
$ws.Range("B12").Value = 85
$ws.Range("C12").Value = -2.4
$ws.Range("E12").Value = "82.6/140"
